# relecture et correction donnees metabo et analyses de survie
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 1: headers -------------------------------------------------
$ws.Range("D1").Value = "label_survie_all"
$ws.Range("E1").Value = "label_survie_selected"
$ws.Range("F1").Value = "label_survie_index"

# --- Column C (label_index) : replace old mtv_* labels with index_* ---
$ws.Range("C2").Value = "index_4_0"
$ws.Range("C3").Value = "index_4_2"
$ws.Range("C4").Value = "index_2_0"
$ws.Range("C5").Value = "index_suv_4_0"
$ws.Range("C6").Value = "index_suv_4_2"
$ws.Range("C7").Value = "index_suv_2_0"

# --- Column E (label_survie_selected), new column -------------------
$ws.Range("E2").Value = "sum_cibles_tep_0_mtv1"
$ws.Range("E3").Value = "sum_cibles_tep_0_mtv2"
$ws.Range("E4").Value = "sum_cibles_tep_2_mtv1"
$ws.Range("E5").Value = "sum_cibles_tep_2_mtv2"
$ws.Range("E6").Value = "sum_cibles_tep_4_mtv1"
$ws.Range("E7").Value = "sum_cibles_tep_4_mtv2"
$ws.Range("E8").Value = "moyenne_cibles_tep_0_suv_peak"
$ws.Range("E9").Value = "moyenne_cibles_tep_2_suv_peak"
$ws.Range("E10").Value = "moyenne_cibles_tep_4_suv_peak"

# --- Column F (label_survie_index), new column -----------------------
$ws.Range("F2").Value = "index_4_0"
$ws.Range("F3").Value = "index_4_2"
$ws.Range("F4").Value = "index_2_0"
$ws.Range("F5").Value = "index_suv_4_0"
$ws.Range("F6").Value = "index_suv_4_2"
$ws.Range("F7").Value = "index_suv_2_0"

# --- Column D (label_survie_all) gains the new index_* rows at the end
$ws.Range("D26").Value = "index_4_0"
$ws.Range("D27").Value = "index_4_2"
$ws.Range("D28").Value = "index_2_0"
$ws.Range("D29").Value = "index_suv_4_0"
$ws.Range("D30").Value = "index_suv_4_2"
$ws.Range("D31").Value = "index_suv_2_0"

# --- column widths: D widened, new column E sized -------------------
$ws.Columns.Item(4).ColumnWidth = 34.5
$ws.Columns.Item(5).ColumnWidth = 22.166666666666668

# --- selection moved to C12 ------------------------------------------
$ws.Range("C12").Select()
